# Update NATMI TPM-derived metrics (Ccl11-Ccr3) with refreshed computed values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 2.940931333333333
    "H2" = 8.822794
    "I2" = 0.03454096854573427
    "J2" = 0.03454096854573427
    "M2" = 0.1808983333333333
    "N2" = 0.5426949999999999
    "O2" = 0.09546831801815302
    "P2" = 0.09546831801815302
    "Q2" = 0.5320095766477777
    "R2" = 4.78808618983
    "S2" = 0.003297568169779179
    "T2" = 0.003297568169779179
    "G3" = 2.940931333333333
    "H3" = 8.822794
    "I3" = 0.03454096854573427
    "J3" = 0.03454096854573427
    "M3" = 1.572737
    "N3" = 4.718211
    "O3" = 0.8300051930177132
    "P3" = 0.8300051930177132
    "Q3" = 4.625311522392667
    "R3" = 41.627803701534
    "S3" = 0.02866918326482093
    "T3" = 0.02866918326482093
    "G4" = 2.940931333333333
    "H4" = 8.822794
    "I4" = 0.03454096854573427
    "J4" = 0.03454096854573427
    "K4" = 1
    "L4" = 0.3333333333333333
    "M4" = 0.1412166666666667
    "N4" = 0.42365
    "O4" = 0.07452648896413371
    "P4" = 0.07452648896413369
    "Q4" = 0.4153085197888889
    "R4" = 3.7377766781
    "S4" = 0.002574217111134155
    "T4" = 0.002574217111134154
    "I5" = 0.8984069237831173
    "J5" = 0.8984069237831174
    "M5" = 0.1808983333333333
    "N5" = 0.5426949999999999
    "O5" = 0.09546831801815302
    "P5" = 0.09546831801815302
    "Q5" = 13.83751259164722
    "R5" = 124.537613324825
    "S5" = 0.0857693979094372
    "T5" = 0.0857693979094372
    "I6" = 0.8984069237831173
    "J6" = 0.8984069237831174
    "M6" = 1.572737
    "N6" = 4.718211
    "O6" = 0.8300051930177132
    "P6" = 0.8300051930177132
    "Q6" = 120.3038615106983
    "R6" = 1082.734753596285
    "S6" = 0.7456824121830562
    "T6" = 0.7456824121830563
    "I7" = 0.8984069237831173
    "J7" = 0.8984069237831174
    "K7" = 1
    "L7" = 0.3333333333333333
    "M7" = 0.1412166666666667
    "N7" = 0.42365
    "O7" = 0.07452648896413371
    "P7" = 0.07452648896413369
    "Q7" = 10.80213049586111
    "R7" = 97.21917446275
    "S7" = 0.0669551136906238
    "T7" = 0.06695511369062379
    "G8" = 1.002166333333333
    "H8" = 3.006499
    "I8" = 0.01177035159063915
    "J8" = 0.01177035159063915
    "M8" = 0.1808983333333333
    "N8" = 0.5426949999999999
    "O8" = 0.09546831801815302
    "P8" = 0.09546831801815302
    "Q8" = 0.1812902194227777
    "R8" = 1.631611974805
    "S8" = 0.001123695668840611
    "T8" = 0.001123695668840611
    "G9" = 1.002166333333333
    "H9" = 3.006499
    "I9" = 0.01177035159063915
    "J9" = 0.01177035159063915
    "M9" = 1.572737
    "N9" = 4.718211
    "O9" = 0.8300051930177132
    "P9" = 0.8300051930177132
    "Q9" = 1.576144072587667
    "R9" = 14.185296653289
    "S9" = 0.009769452943874794
    "T9" = 0.009769452943874796
    "G10" = 1.002166333333333
    "H10" = 3.006499
    "I10" = 0.01177035159063915
    "J10" = 0.01177035159063915
    "K10" = 1
    "L10" = 0.3333333333333333
    "M10" = 0.1412166666666667
    "N10" = 0.42365
    "O10" = 0.07452648896413371
    "P10" = 0.07452648896413369
    "Q10" = 0.1415225890388889
    "R10" = 1.27370330135
    "S10" = 0.000877202977923742
    "T10" = 0.000877202977923742
    "G11" = 2.356521666666667
    "H11" = 7.069565
    "I11" = 0.02767713065691252
    "J11" = 0.02767713065691253
    "M11" = 0.1808983333333333
    "N11" = 0.5426949999999999
    "O11" = 0.09546831801815302
    "P11" = 0.09546831801815302
    "Q11" = 0.4262908419638888
    "R11" = 3.836617577674999
    "S11" = 0.002642289111384097
    "T11" = 0.002642289111384097
    "G12" = 2.356521666666667
    "H12" = 7.069565
    "I12" = 0.02767713065691252
    "J12" = 0.02767713065691253
    "M12" = 1.572737
    "N12" = 4.718211
    "O12" = 0.8300051930177132
    "P12" = 0.8300051930177132
    "Q12" = 3.706188816468333
    "R12" = 33.355699348215
    "S12" = 0.02297216217306715
    "T12" = 0.02297216217306715
    "G13" = 2.356521666666667
    "H13" = 7.069565
    "I13" = 0.02767713065691252
    "J13" = 0.02767713065691253
    "K13" = 1
    "L13" = 0.3333333333333333
    "M13" = 0.1412166666666667
    "N13" = 0.42365
    "O13" = 0.07452648896413371
    "P13" = 0.07452648896413369
    "Q13" = 0.3327801346944445
    "R13" = 2.99502121225
    "S13" = 0.002062679372461278
    "T13" = 0.002062679372461278
    "E14" = 2
    "F14" = 0.6666666666666666
    "G14" = 0.5156633333333334
    "H14" = 1.54699
    "I14" = 0.006056418514425867
    "J14" = 0.006056418514425868
    "M14" = 0.1808983333333333
    "N14" = 0.5426949999999999
    "O14" = 0.09546831801815302
    "P14" = 0.09546831801815302
    "Q14" = 0.0932826375611111
    "R14" = 0.8395437380499999
    "S14" = 0.0005781960887862386
    "T14" = 0.0005781960887862387
    "E15" = 2
    "F15" = 0.6666666666666666
    "G15" = 0.5156633333333334
    "H15" = 1.54699
    "I15" = 0.006056418514425867
    "J15" = 0.006056418514425868
    "M15" = 1.572737
    "N15" = 4.718211
    "O15" = 0.8300051930177132
    "P15" = 0.8300051930177132
    "Q15" = 0.8110028038766668
    "R15" = 7.29902523489
    "S15" = 0.005026858818062094
    "T15" = 0.005026858818062095
    "E16" = 2
    "F16" = 0.6666666666666666
    "G16" = 0.5156633333333334
    "H16" = 1.54699
    "I16" = 0.006056418514425867
    "J16" = 0.006056418514425868
    "K16" = 1
    "L16" = 0.3333333333333333
    "M16" = 0.1412166666666667
    "N16" = 0.42365
    "O16" = 0.07452648896413371
    "P16" = 0.07452648896413369
    "Q16" = 0.07282025705555557
    "R16" = 0.6553823135000001
    "S16" = 0.0004513636075775344
    "T16" = 0.0004513636075775344
    "E17" = 3
    "F17" = 1
    "G17" = 1.834685
    "H17" = 5.504055
    "I17" = 0.02154820690917088
    "J17" = 0.02154820690917089
    "M17" = 0.1808983333333333
    "N17" = 0.5426949999999999
    "O17" = 0.09546831801815302
    "P17" = 0.09546831801815302
    "Q17" = 0.3318914586916666
    "R17" = 2.987023128225
    "S17" = 0.002057171069925688
    "T17" = 0.002057171069925689
    "E18" = 3
    "F18" = 1
    "G18" = 1.834685
    "H18" = 5.504055
    "I18" = 0.02154820690917088
    "J18" = 0.02154820690917089
    "M18" = 1.572737
    "N18" = 4.718211
    "O18" = 0.8300051930177132
    "P18" = 0.8300051930177132
    "Q18" = 2.885476982845
    "R18" = 25.969292845605
    "S18" = 0.017885123634832
    "T18" = 0.01788512363483201
    "E19" = 3
    "F19" = 1
    "G19" = 1.834685
    "H19" = 5.504055
    "I19" = 0.02154820690917088
    "J19" = 0.02154820690917089
    "K19" = 1
    "L19" = 0.3333333333333333
    "M19" = 0.1412166666666667
    "N19" = 0.42365
    "O19" = 0.07452648896413371
    "P19" = 0.07452648896413369
    "Q19" = 0.2590881000833334
    "R19" = 2.33179290075
    "S19" = 0.001605912204413194
    "T19" = 0.001605912204413194
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
